# Apply country-data refresh (new countries added, case counts updated as of 31-Mar-2020 11:50)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 11:50"

$ws.Range("D13").Value = 1696
$ws.Range("E13").Value = 10374

$ws.Range("D28").Value = 537
$ws.Range("E28").Value = 2186

$ws.Range("A30").Value = "Rusia"
$ws.Range("B30").Value = 2337
$ws.Range("C30").Value = 501
$ws.Range("D30").Value = 121
$ws.Range("E30").Value = 2199
$ws.Range("F30").Value = 8
$ws.Range("G30").Value = 8
$ws.Range("H30").Value = 17

$ws.Range("A31").Value = "Polonia"
$ws.Range("B31").Value = 2132
$ws.Range("C31").Value = 77
$ws.Range("D31").Value = 7
$ws.Range("E31").Value = 2094
$ws.Range("F31").Value = 3
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 31

$ws.Range("A32").Value = "Rumania"
$ws.Range("B32").Value = 2109
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 209
$ws.Range("E32").Value = 1831
$ws.Range("F32").Value = 33
$ws.Range("G32").Value = 4
$ws.Range("H32").Value = 69

$ws.Range("A33").Value = "Filipinas"
$ws.Range("B33").Value = 2084
$ws.Range("C33").Value = 538
$ws.Range("D33").Value = 49
$ws.Range("E33").Value = 1947
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 10
$ws.Range("H33").Value = 88

$ws.Range("A34").Value = "Luxemburgo"
$ws.Range("B34").Value = 1988
$ws.Range("D34").Value = 80
$ws.Range("E34").Value = 1886
$ws.Range("F34").Value = 31
$ws.Range("H34").Value = 22

$ws.Range("A35").Value = "Ecuador"
$ws.Range("B35").Value = 1966
$ws.Range("D35").Value = 54
$ws.Range("E35").Value = 1850
$ws.Range("F35").Value = 58
$ws.Range("H35").Value = 62

$ws.Range("A36").Value = "Japon"
$ws.Range("B36").Value = 1953
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 424
$ws.Range("E36").Value = 1473
$ws.Range("F36").Value = 56
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 56

$ws.Range("A37").Value = "Pakistan"
$ws.Range("B37").Value = 1865
$ws.Range("C37").Value = 148
$ws.Range("D37").Value = 76
$ws.Range("E37").Value = 1764
$ws.Range("F37").Value = 12
$ws.Range("G37").Value = 4
$ws.Range("H37").Value = 25

$ws.Range("A39").Value = "Indonesia"
$ws.Range("B39").Value = 1528
$ws.Range("C39").Value = 114
$ws.Range("D39").Value = 81
$ws.Range("E39").Value = 1311
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 14
$ws.Range("H39").Value = 136

$ws.Range("A40").Value = "Arabia Saudita"
$ws.Range("B40").Value = 1453
$ws.Range("D40").Value = 115
$ws.Range("E40").Value = 1330
$ws.Range("F40").Value = 12
$ws.Range("H40").Value = 8

$ws.Range("A68").Value = "Armenia"
$ws.Range("B68").Value = 532
$ws.Range("C68").Value = 50
$ws.Range("D68").Value = 30
$ws.Range("E68").Value = 499
$ws.Range("F68").Value = 15
$ws.Range("H68").Value = 3

$ws.Range("A69").Value = "Barein"
$ws.Range("B69").Value = 515
$ws.Range("C69").Value = 0
$ws.Range("D69").Value = 295
$ws.Range("E69").Value = 216
$ws.Range("F69").Value = 2
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 4

$ws.Range("A70").Value = "Hungria"
$ws.Range("B70").Value = 492
$ws.Range("C70").Value = 45
$ws.Range("D70").Value = 37
$ws.Range("E70").Value = 439
$ws.Range("F70").Value = 6
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 16

$ws.Range("B71").Value = 463
$ws.Range("C71").Value = 17
$ws.Range("E71").Value = 416
$ws.Range("F71").Value = 7

$ws.Range("A78").Value = "Kazajistan"
$ws.Range("B78").Value = 336
$ws.Range("C78").Value = 34
$ws.Range("D78").Value = 22
$ws.Range("E78").Value = 312
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 1

$ws.Range("A79").Value = "Costa Rica"
$ws.Range("B79").Value = 330
$ws.Range("C79").Value = 0
$ws.Range("D79").Value = 4
$ws.Range("E79").Value = 324
$ws.Range("F79").Value = 7
$ws.Range("G79").Value = 0

$ws.Range("D82").Value = 18
$ws.Range("E82").Value = 278

$ws.Range("A88").Value = "Albania"
$ws.Range("B88").Value = 243
$ws.Range("C88").Value = 20
$ws.Range("D88").Value = 52
$ws.Range("E88").Value = 178
$ws.Range("F88").Value = 8
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = 13

$ws.Range("A89").Value = "Republica de Chipre"
$ws.Range("D89").Value = 22
$ws.Range("E89").Value = 201
$ws.Range("F89").Value = 3
$ws.Range("H89").Value = 7

$ws.Range("A90").Value = "San Marino"
$ws.Range("B90").Value = 230
$ws.Range("D90").Value = 13
$ws.Range("E90").Value = 192
$ws.Range("F90").Value = 16
$ws.Range("H90").Value = 25

$ws.Range("A91").Value = "Reunion"
$ws.Range("B91").Value = 224
$ws.Range("D91").Value = 1
$ws.Range("E91").Value = 223
$ws.Range("F91").Value = 0
$ws.Range("H91").Value = 0

$ws.Range("B108").Value = 129
$ws.Range("C108").Value = 2
$ws.Range("D108").Value = 45
$ws.Range("E108").Value = 83

$ws.Range("F146").Value = 2

$ws.Range("A153").Value = "Nueva Caledonia"
$ws.Range("C153").Value = 1
$ws.Range("E153").Value = 16
$ws.Range("H153").Value = 0

$ws.Range("A154").Value = "Gabon"
$ws.Range("B154").Value = 16
$ws.Range("C154").Value = 9
$ws.Range("H154").Value = 1

$ws.Range("A155").Value = "Eritrea"

$ws.Range("A179").Value = "Antigua y Barbuda"
$ws.Range("C179").Value = 0

$ws.Range("A181").Value = "Republica del Chad"
$ws.Range("C181").Value = 2

$ws.Range("A182").Value = "Sudan"
$ws.Range("C182").Value = 1
$ws.Range("D182").Value = 0
$ws.Range("E182").Value = 5

$ws.Range("A183").Value = "Angola"
$ws.Range("B183").Value = 7
$ws.Range("D183").Value = 1
$ws.Range("E183").Value = 4
$ws.Range("H183").Value = 2

$ws.Range("A184").Value = "San Martin (Parte Holandesa)"

$ws.Range("A185").Value = "Santa Sede"
$ws.Range("D185").Value = 0
$ws.Range("E185").Value = 6

$ws.Range("A186").Value = "Benin"
$ws.Range("D186").Value = 1
$ws.Range("H186").Value = 0

$ws.Range("A187").Value = "Cabo Verde"
$ws.Range("D187").Value = 0
$ws.Range("H187").Value = 1

$ws.Range("A188").Value = "San Bartolome"
$ws.Range("D188").Value = 1
$ws.Range("E188").Value = 5
$ws.Range("H188").Value = 0

$ws.Range("A190").Value = "Montserrat"

$ws.Range("A192").Value = "Islas Turcas y Caicos"

$ws.Range("A198").Value = "Belice"

$ws.Range("A199").Value = "Botsuana"

$ws.Range("A200").Value = "Liberia"
$ws.Range("C200").Value = 0

$ws.Range("A201").Value = "Islas Virgenes Britanicas"
$ws.Range("C201").Value = 1
